$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of user data (row 3), mirroring the inline-string layout used by
# rows 1-2: Name, Username, PasswordHash, Phone, City, PostalCode.
$ws.Range("A3").Value = "dadadada"
$ws.Range("B3").Value = "dadadsa"
$ws.Range("C3").Value = "a03ab19b866fc585b5cb1812a2f63ca861e7e7643ee5d43fd7106b623725fd67"

# D3/F3 look numeric ("123344444" / "12312") but must stay text, like the
# rest of the sheet (t="inlineStr" in the target) - force text formatting
# before assigning so Excel doesn't coerce them into numbers.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "123344444"

$ws.Range("E3").Value = "Gonçalo"

$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "12312"
